$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 10375
$ws.Range("I32").Value = 8500
$ws.Range("J32").Value = 12250
$ws.Range("K32").Value = 8500
$ws.Range("L32").Value = 12250
$ws.Range("M32").Value = -8174
$ws.Range("N32").Value = -12902

$ws.Range("H64").Value = 4545.4546
$ws.Range("I64").Value = 7000
$ws.Range("J64").Value = 4000
$ws.Range("K64").Value = 7000
$ws.Range("L64").Value = 4000
$ws.Range("M64").Value = -6752
$ws.Range("N64").Value = -4496

$ws.Range("H67").Value = 4545.4546
$ws.Range("I67").Value = 7000
$ws.Range("J67").Value = 4000
$ws.Range("K67").Value = 7000
$ws.Range("L67").Value = 4000
$ws.Range("M67").Value = -6142
$ws.Range("N67").Value = -5716

$ws.Range("H70").Value = 999
$ws.Range("I70").Value = 999
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 2997
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -2727
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 999
$ws.Range("I73").Value = 999
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 2997
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -2061
$ws.Range("N73").ClearContents()

$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("M76").ClearContents()
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("N79").ClearContents()

$ws.Range("H100").Value = 1885.4445
$ws.Range("I100").Value = 2660
$ws.Range("J100").Value = 1498.1666
$ws.Range("K100").Value = 2660
$ws.Range("L100").Value = 1498.1666
$ws.Range("M100").Value = -2119
$ws.Range("N100").Value = -2580.1666

$ws.Range("H107").Value = 306.16666
$ws.Range("I107").Value = 327.6
$ws.Range("J107").Value = 199
$ws.Range("K107").Value = 327.6
$ws.Range("L107").Value = 199
$ws.Range("M107").Value = 1592.4
$ws.Range("N107").Value = -4039

$ws.Range("H112").Value = 1643.8572
$ws.Range("I112").Value = 766.6667
$ws.Range("J112").Value = 1790.0555
$ws.Range("K112").Value = 2300.0001
$ws.Range("L112").Value = 5370.166499999999
$ws.Range("M112").Value = -1192.0001
$ws.Range("N112").Value = -7586.166499999999

$ws.Range("H135").Value = 2200
$ws.Range("I135").Value = 2200
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 19800
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -17265

$ws.Range("H138").Value = 2395.111
$ws.Range("I138").Value = 1897.8334
$ws.Range("J138").Value = 2643.75
$ws.Range("K138").Value = 5693.5002
$ws.Range("L138").Value = 7931.25
$ws.Range("M138").Value = -553.5002000000004
$ws.Range("N138").Value = -18211.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1259.4615
$ws.Range("I2").Value = 1326.2727
$ws.Range("J2").Value = 892
$ws.Range("K2").Value = 1326.2727
$ws.Range("L2").Value = 892
$ws.Range("M2").Value = -1213.2727
$ws.Range("N2").Value = -1118

$ws.Range("H61").Value = 3248.5
$ws.Range("I61").Value = 2498.25
$ws.Range("J61").Value = 3998.75
$ws.Range("K61").Value = 2498.25
$ws.Range("L61").Value = 3998.75
$ws.Range("M61").Value = -2286.25
$ws.Range("N61").Value = -4422.75

$ws.Range("H74").Value = 22215660
$ws.Range("I74").Value = 28561134
$ws.Range("J74").Value = 6499.5
$ws.Range("K74").Value = 28561134
$ws.Range("L74").Value = 6499.5
$ws.Range("M74").Value = -28560260
$ws.Range("N74").Value = -8247.5

$ws.Range("H77").Value = 22215660
$ws.Range("I77").Value = 28561134
$ws.Range("J77").Value = 6499.5
$ws.Range("K77").Value = 142805670
$ws.Range("L77").Value = 32497.5
$ws.Range("M77").Value = -142801302
$ws.Range("N77").Value = -41233.5

$ws.Range("H110").Value = 1255.25
$ws.Range("I110").Value = 1255.25
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 1255.25
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 789.75

$ws.Range("H116").Value = 1259.4615
$ws.Range("I116").Value = 1326.2727
$ws.Range("J116").Value = 892
$ws.Range("K116").Value = 1326.2727
$ws.Range("L116").Value = 892
$ws.Range("M116").Value = 967.7273
$ws.Range("N116").Value = -5480

$ws.Range("H136").Value = 3248.5
$ws.Range("I136").Value = 2498.25
$ws.Range("J136").Value = 3998.75
$ws.Range("K136").Value = 7494.75
$ws.Range("L136").Value = 11996.25
$ws.Range("M136").Value = -4944.75
$ws.Range("N136").Value = -17096.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1259.4615
$ws.Range("I3").Value = 1326.2727
$ws.Range("J3").Value = 892
$ws.Range("K3").Value = 1326.2727
$ws.Range("L3").Value = 892
$ws.Range("M3").Value = -1212.2727
$ws.Range("N3").Value = -1120

$ws.Range("H105").Value = 3626.818
$ws.Range("I105").Value = 3411.2222
$ws.Range("J105").Value = 4597
$ws.Range("K105").Value = 3411.2222
$ws.Range("L105").Value = 4597
$ws.Range("M105").Value = -1664.2222
$ws.Range("N105").Value = -8091

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5998.6
$ws.Range("I99").Value = 5998
$ws.Range("J99").Value = 5998.75
$ws.Range("K99").Value = 5998
$ws.Range("L99").Value = 5998.75
$ws.Range("M99").Value = -4500
$ws.Range("N99").Value = -8994.75

$ws.Range("H107").Value = 1937.8182
$ws.Range("I107").Value = 1120
$ws.Range("J107").Value = 2405.1428
$ws.Range("K107").Value = 1120
$ws.Range("L107").Value = 2405.1428
$ws.Range("M107").Value = 800
$ws.Range("N107").Value = -6245.1428

$ws.Range("H126").Value = 5998.6
$ws.Range("I126").Value = 5998
$ws.Range("J126").Value = 5998.75
$ws.Range("K126").Value = 17994
$ws.Range("L126").Value = 17996.25
$ws.Range("M126").Value = -15524
$ws.Range("N126").Value = -22936.25

$ws.Range("H141").Value = 53633.43
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 53633.43
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 53633.43
$ws.Range("N141").Value = -63993.43

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 4948.5
$ws.Range("I81").Value = 4900
$ws.Range("J81").Value = 4997
$ws.Range("K81").Value = 14700
$ws.Range("L81").Value = 14991
$ws.Range("M81").Value = -13577
$ws.Range("N81").Value = -17237

$ws.Range("H84").Value = 4948.5
$ws.Range("I84").Value = 4900
$ws.Range("J84").Value = 4997
$ws.Range("K84").Value = 44100
$ws.Range("L84").Value = 44973
$ws.Range("M84").Value = -38484
$ws.Range("N84").Value = -56205

$ws.Range("H131").Value = 1263.25
$ws.Range("I131").Value = 851.6667
$ws.Range("J131").Value = 2498
$ws.Range("K131").Value = 2555.0001
$ws.Range("L131").Value = 7494
$ws.Range("M131").Value = 2484.9999
$ws.Range("N131").Value = -17574

$ws.Range("H140").Value = 1304.6
$ws.Range("I140").Value = 1304.6
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 3913.8
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 1266.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3250
$ws.Range("I122").Value = 2333.3333
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 6999.999899999999
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -4549.999899999999
$ws.Range("N122").Value = -22900

$ws.Range("H132").Value = 2572.05
$ws.Range("I132").Value = 1703.9166
$ws.Range("J132").Value = 3874.25
$ws.Range("K132").Value = 5111.7498
$ws.Range("L132").Value = 11622.75
$ws.Range("M132").Value = -2581.7498
$ws.Range("N132").Value = -16682.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3466.6667
$ws.Range("I7").Value = 3466.6667
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3466.6667
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -3354.6667

$ws.Range("H122").Value = 6767.3335
$ws.Range("I122").Value = 6738.375
$ws.Range("J122").Value = 6999
$ws.Range("K122").Value = 20215.125
$ws.Range("L122").Value = 20997
$ws.Range("M122").Value = -17765.125
$ws.Range("N122").Value = -25897

$ws.Range("H126").Value = 3466.6667
$ws.Range("I126").Value = 3466.6667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10400.0001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -7930.000100000001

$ws.Range("H132").Value = 4108.9
$ws.Range("I132").Value = 3458.6
$ws.Range("J132").Value = 4759.2
$ws.Range("K132").Value = 10375.8
$ws.Range("L132").Value = 14277.6
$ws.Range("M132").Value = -7845.799999999999
$ws.Range("N132").Value = -19337.6

$ws.Range("H136").Value = 16000799
$ws.Range("I136").Value = 16000799
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 48002397
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -47999847

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 30000.5
$ws.Range("I4").Value = 28334
$ws.Range("J4").Value = 35000
$ws.Range("K4").Value = 28334
$ws.Range("L4").Value = 35000
$ws.Range("M4").Value = -28221
$ws.Range("N4").Value = -35226

$ws.Range("H62").Value = 134997.5
$ws.Range("I62").Value = 173330
$ws.Range("J62").Value = 20000
$ws.Range("K62").Value = 173330
$ws.Range("L62").Value = 20000
$ws.Range("M62").Value = -172706
$ws.Range("N62").Value = -21248

$ws.Range("H65").Value = 134997.5
$ws.Range("I65").Value = 173330
$ws.Range("J65").Value = 20000
$ws.Range("K65").Value = 866650
$ws.Range("L65").Value = 100000
$ws.Range("M65").Value = -863530
$ws.Range("N65").Value = -106240

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H81").Value = 12582.5
$ws.Range("I81").Value = 9375.25
$ws.Range("J81").Value = 18997
$ws.Range("K81").Value = 18750.5
$ws.Range("L81").Value = 37994
$ws.Range("M81").Value = -17689.5
$ws.Range("N81").Value = -40116

$ws.Range("H84").Value = 12582.5
$ws.Range("I84").Value = 9375.25
$ws.Range("J84").Value = 18997
$ws.Range("K84").Value = 93752.5
$ws.Range("L84").Value = 189970
$ws.Range("M84").Value = -88448.5
$ws.Range("N84").Value = -200578

$ws.Range("H136").Value = 2030.1724
$ws.Range("I136").Value = 1736.6957
$ws.Range("J136").Value = 3155.1667
$ws.Range("K136").Value = 5210.0871
$ws.Range("L136").Value = 9465.500100000001
$ws.Range("M136").Value = -2660.0871
$ws.Range("N136").Value = -14565.5001
